# Scheduled runner update: refresh market-board derived profit figures
# (currentAveragePrice / LevePrice / LeveProfit columns) across the
# per-job Leve sheets, matching the latest price pull.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 805.06665
$ws.Range("I38").Value = 407.9
$ws.Range("K38").Value = 1223.7
$ws.Range("M38").Value = -851.6999999999998
$ws.Range("H58").Value = 2003.5714
$ws.Range("I58").Value = 63
$ws.Range("J58").Value = 2779.8
$ws.Range("K58").Value = 189
$ws.Range("L58").Value = 8339.400000000001
$ws.Range("M58").Value = -39
$ws.Range("N58").Value = -8639.400000000001
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1069.6
$ws.Range("I2").Value = 962
$ws.Range("J2").Value = 1500
$ws.Range("K2").Value = 962
$ws.Range("L2").Value = 1500
$ws.Range("M2").Value = -849
$ws.Range("N2").Value = -1726
$ws.Range("H61").Value = 3849.5
$ws.Range("I61").Value = 3799.6667
$ws.Range("K61").Value = 3799.6667
$ws.Range("M61").Value = -3587.6667
$ws.Range("H116").Value = 1069.6
$ws.Range("I116").Value = 962
$ws.Range("J116").Value = 1500
$ws.Range("K116").Value = 962
$ws.Range("L116").Value = 1500
$ws.Range("M116").Value = 1332
$ws.Range("N116").Value = -6088
$ws.Range("H124").Value = 76999
$ws.Range("J124").Value = 76999
$ws.Range("L124").Value = 76999
$ws.Range("N124").Value = -86819
$ws.Range("H136").Value = 3849.5
$ws.Range("I136").Value = 3799.6667
$ws.Range("K136").Value = 11399.0001
$ws.Range("M136").Value = -8849.000100000001
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1069.6
$ws.Range("I3").Value = 962
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 962
$ws.Range("L3").Value = 1500
$ws.Range("M3").Value = -848
$ws.Range("N3").Value = -1728
$ws.Range("H40").Value = 59999
$ws.Range("J40").Value = 59999
$ws.Range("L40").Value = 59999
$ws.Range("N40").Value = -60529
$ws.Range("H80").Value = 482.81818
$ws.Range("I80").Value = 186.16667
$ws.Range("K80").Value = 186.16667
$ws.Range("M80").Value = 811.8333299999999
$ws.Range("H83").Value = 482.81818
$ws.Range("I83").Value = 186.16667
$ws.Range("K83").Value = 930.8333500000001
$ws.Range("M83").Value = 4061.16665
$ws.Range("H134").Value = 1233.3334
$ws.Range("I134").Value = 1233.3334
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 3700.0002
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -1165.0002
$ws.Range("N134").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 920.125
$ws.Range("I31").Value = 938.8570999999999
$ws.Range("J31").Value = 789
$ws.Range("K31").Value = 938.8570999999999
$ws.Range("L31").Value = 789
$ws.Range("M31").Value = -643.8570999999999
$ws.Range("N31").Value = -1379
$ws.Range("H34").Value = 920.125
$ws.Range("I34").Value = 938.8570999999999
$ws.Range("J34").Value = 789
$ws.Range("K34").Value = 938.8570999999999
$ws.Range("L34").Value = 789
$ws.Range("M34").Value = -736.8570999999999
$ws.Range("N34").Value = -1193
$ws.Range("H58").Value = 279.4
$ws.Range("I58").Value = 279.4
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 279.4
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -76.39999999999998
$ws.Range("N58").ClearContents()
$ws.Range("H62").Value = 2720.5
$ws.Range("I62").Value = 2830.75
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2830.75
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -2206.75
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2720.5
$ws.Range("I65").Value = 2830.75
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 14153.75
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -11033.75
$ws.Range("N65").Value = -18740
$ws.Range("H86").Value = 333336500
$ws.Range("I86").Value = 333336500
$ws.Range("K86").Value = 333336500
$ws.Range("M86").Value = -333335377
$ws.Range("H89").Value = 333336500
$ws.Range("I89").Value = 333336500
$ws.Range("K89").Value = 1666682500
$ws.Range("M89").Value = -1666676884
$ws.Range("H132").Value = 2276.3076
$ws.Range("J132").Value = 2437.5
$ws.Range("L132").Value = 7312.5
$ws.Range("N132").Value = -12372.5
$ws.Range("H136").Value = 279.4
$ws.Range("I136").Value = 279.4
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 838.1999999999999
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = 1711.8
$ws.Range("N136").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 1017.5714
$ws.Range("I38").Value = 162.5
$ws.Range("J38").Value = 1359.6
$ws.Range("K38").Value = 487.5
$ws.Range("L38").Value = 4078.8
$ws.Range("M38").Value = -140.5
$ws.Range("N38").Value = -4772.799999999999
$ws.Range("H75").Value = 2403
$ws.Range("I75").Value = 200
$ws.Range("J75").Value = 2770.1667
$ws.Range("K75").Value = 600
$ws.Range("L75").Value = 8310.500100000001
$ws.Range("M75").Value = 398
$ws.Range("N75").Value = -10306.5001
$ws.Range("H78").Value = 2403
$ws.Range("I78").Value = 200
$ws.Range("J78").Value = 2770.1667
$ws.Range("K78").Value = 1800
$ws.Range("L78").Value = 24931.5003
$ws.Range("M78").Value = 3192
$ws.Range("N78").Value = -34915.5003
$ws.Range("H107").Value = 850
$ws.Range("I107").Value = 847
$ws.Range("J107").Value = 850.75
$ws.Range("K107").Value = 2541
$ws.Range("L107").Value = 2552.25
$ws.Range("M107").Value = -621
$ws.Range("N107").Value = -6392.25
$ws.Range("H118").Value = 2858.1667
$ws.Range("I118").Value = 2858.1667
$ws.Range("K118").Value = 8574.500100000001
$ws.Range("M118").Value = -7331.500100000001
$ws.Range("H139").Value = 2829.8333
$ws.Range("I139").Value = 1744.75
$ws.Range("J139").Value = 5000
$ws.Range("K139").Value = 5234.25
$ws.Range("L139").Value = 15000
$ws.Range("M139").Value = -94.25
$ws.Range("N139").Value = -25280
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1002000
$ws.Range("I97").Value = 0
$ws.Range("K97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("H102").Value = 2940.7144
$ws.Range("I102").Value = 3024.25
$ws.Range("K102").Value = 3024.25
$ws.Range("M102").Value = -1402.25
$ws.Range("H126").Value = 3499
$ws.Range("I126").Value = 3499
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 10497
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -8027
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 42514.145
$ws.Range("J136").Value = 42514.145
$ws.Range("L136").Value = 127542.435
$ws.Range("N136").Value = -132642.435
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H122").Value = 2000
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 1500
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 4500
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -9400
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 500
$ws.Range("I100").Value = 500
$ws.Range("K100").Value = 1000
$ws.Range("M100").Value = -459
